$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$values = @(
    @(1, 3),
    @(3, 5),
    @(4, 7),
    @(5, 6),
    @(5, 5),
    @(5, 6),
    @(2, 5),
    @(3, 5),
    @(6, 8),
    @(8, 8),
    @(2, 4),
    @(4, 5),
    @(9, 9),
    @(8, 9),
    @(6, 7),
    @(6, 8),
    @(5, 6),
    @(6, 7)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
